$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "@prefix": add new prefix rows for pixels / bindata /
# tagAnnotation (and keep unitLength / image together, with unitLength
# now preceding image).
# ------------------------------------------------------------------
$wsPrefix = $wb.Worksheets.Item("@prefix")

# Insert a new row above the existing "image" row (row 13). This
# pushes the old row 13 ("image") down to row 14 and the old row 14
# ("unitLength") down to row 15.
$wsPrefix.Rows.Item(13).Insert() | Out-Null

$wsPrefix.Range("A13").Value = "unitLength"
$wsPrefix.Range("B13").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/Unit/UnitLength#"

# Row 14 already contains the correct "image" prefix values (shifted
# down automatically by the insert above), so nothing else to do there.

# Row 15 currently holds a duplicate of the old "unitLength" row that
# was pushed down by the insert; overwrite it with the new "pixels"
# prefix entry.
$wsPrefix.Range("A15").Value = "pixels"
$wsPrefix.Range("B15").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/tagannotation/pixels/"

# New rows for bindata and tagAnnotation prefixes.
$wsPrefix.Range("A16").Value = "bindata"
$wsPrefix.Range("B16").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/tagannotation/bindata/"

$wsPrefix.Range("A17").Value = "tagAnnotation"
$wsPrefix.Range("B17").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/tagannotation/tagAnnotation/"

# The sheet also gained an explicit page setup (A4, portrait).
$wsPrefix.PageSetup.PaperSize = 9
$wsPrefix.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# Fix the ad-hoc "[prefix:local]" style identifiers so that they use
# the proper "prefix:local" naming convention everywhere they are
# referenced across the workbook.
# ------------------------------------------------------------------

$wsImage = $wb.Worksheets.Item("Image")
$wsImage.Range("E5").Value = "pixels:pixels0.0"
$wsImage.Range("G5").Value = "tagAnnotation:tagAnnotation3"

$wsPixels = $wb.Worksheets.Item("Pixels")
$wsPixels.Range("B5").Value = "pixels:pixels0.0"
$wsPixels.Range("M5").Value = "bindata:bindata1"

$wsBinData = $wb.Worksheets.Item("Binary_Data")
$wsBinData.Range("B5").Value = "bindata:bindata1"

$wsStructAnno = $wb.Worksheets.Item("Structured_Annotations")
$wsStructAnno.Range("C5").Value = "tagAnnotation:tagAnnotation1"
$wsStructAnno.Range("C6").Value = "tagAnnotation:tagAnnotation2"
$wsStructAnno.Range("C7").Value = "tagAnnotation:tagAnnotation3"

$wsTagAnno = $wb.Worksheets.Item("Tag_Annotation")
$wsTagAnno.Range("B5").Value = "tagAnnotation:tagAnnotation1"
$wsTagAnno.Range("B6").Value = "tagAnnotation:tagAnnotation2"
$wsTagAnno.Range("B7").Value = "tagAnnotation:tagAnnotation3"
$wsTagAnno.Range("B8").Value = "tagAnnotation:tagAnnotation3"
$wsTagAnno.Range("G7").Value = "tagAnnotation:tagAnnotation1"
$wsTagAnno.Range("G8").Value = "tagAnnotation:tagAnnotation2"
